$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.1707317073170732
$ws.Range("C2").Value2 = 0.6189024390243902
$ws.Range("J2").Value2 = 0.02439024390243903
$ws.Range("P2").Value2 = 0.1310975609756098
$ws.Range("S2").Value2 = 0.05487804878048781
$ws.Range("B3").Value2 = 0.004739336492890996
$ws.Range("C3").Value2 = 0.01421800947867299
$ws.Range("J3").Value2 = 0.04265402843601896
$ws.Range("P3").Value2 = 0.8056872037914692
$ws.Range("S3").Value2 = 0.1327014218009479
$ws.Range("J4").Value2 = 0.075
$ws.Range("P4").Value2 = 0.725
$ws.Range("S4").Value2 = 0.2
$ws.Range("P5").Value2 = 1
$ws.Range("B6").Value2 = 0.04313725490196078
$ws.Range("D6").Value2 = 0.02745098039215686
$ws.Range("F6").Value2 = 0.04705882352941176
$ws.Range("J6").Value2 = 0.2705882352941176
$ws.Range("O6").Value2 = 0.01568627450980392
$ws.Range("Q6").Value2 = 0.196078431372549
$ws.Range("R6").Value2 = 0.09803921568627451
$ws.Range("S6").Value2 = 0.3019607843137255
$ws.Range("B7").Value2 = 0.06976744186046512
$ws.Range("D7").Value2 = 0.05581395348837209
$ws.Range("F7").Value2 = 0.07441860465116279
$ws.Range("J7").Value2 = 0.1488372093023256
$ws.Range("O7").Value2 = 0.02325581395348837
$ws.Range("Q7").Value2 = 0.2186046511627907
$ws.Range("R7").Value2 = 0.07906976744186046
$ws.Range("S7").Value2 = 0.3302325581395349
$ws.Range("B8").Value2 = 0.08413001912045889
$ws.Range("D8").Value2 = 0.02294455066921606
$ws.Range("F8").Value2 = 0.08221797323135756
$ws.Range("J8").Value2 = 0.1338432122370937
$ws.Range("O8").Value2 = 0.02294455066921606
$ws.Range("Q8").Value2 = 0.1701720841300191
$ws.Range("R8").Value2 = 0.1089866156787763
$ws.Range("S8").Value2 = 0.3747609942638623
$ws.Range("B9").Value2 = 0.09278350515463918
$ws.Range("D9").Value2 = 0.02577319587628866
$ws.Range("E9").Value2 = 0.005154639175257732
$ws.Range("F9").Value2 = 0.08762886597938144
$ws.Range("J9").Value2 = 0.134020618556701
$ws.Range("O9").Value2 = 0.0154639175257732
$ws.Range("Q9").Value2 = 0.2525773195876289
$ws.Range("R9").Value2 = 0.05670103092783505
$ws.Range("S9").Value2 = 0.3298969072164948
$ws.Range("B10").Value2 = 0.1038575667655786
$ws.Range("D10").Value2 = 0.02611275964391692
$ws.Range("F10").Value2 = 0.05934718100890208
$ws.Range("J10").Value2 = 0.1258160237388724
$ws.Range("O10").Value2 = 0.01068249258160237
$ws.Range("Q10").Value2 = 0.2249258160237389
$ws.Range("R10").Value2 = 0.1032640949554896
$ws.Range("S10").Value2 = 0.3459940652818991
$ws.Range("G11").Value2 = 0.1437699680511182
$ws.Range("J11").Value2 = 0.08626198083067092
$ws.Range("K11").Value2 = 0.207667731629393
$ws.Range("L11").Value2 = 0.5527156549520766
$ws.Range("S11").Value2 = 0.009584664536741214
$ws.Range("G12").Value2 = 0.7388888888888889
$ws.Range("J12").Value2 = 0.1722222222222222
$ws.Range("K12").Value2 = 0.005555555555555556
$ws.Range("L12").Value2 = 0.05
$ws.Range("S12").Value2 = 0.03333333333333333
$ws.Range("F15").Value2 = 0.02135231316725979
$ws.Range("H15").Value2 = 0.1565836298932384
$ws.Range("I15").Value2 = 0.05693950177935943
$ws.Range("J15").Value2 = 0.4448398576512456
$ws.Range("K15").Value2 = 0.0498220640569395
$ws.Range("M15").Value2 = 0.007117437722419928
$ws.Range("O15").Value2 = 0.05693950177935943
$ws.Range("S15").Value2 = 0.2064056939501779
$ws.Range("F16").Value2 = 0.03162055335968379
$ws.Range("H16").Value2 = 0.1897233201581028
$ws.Range("I16").Value2 = 0.06719367588932806
$ws.Range("J16").Value2 = 0.458498023715415
$ws.Range("K16").Value2 = 0.08695652173913043
$ws.Range("M16").Value2 = 0.02371541501976284
$ws.Range("N16").Value2 = 0.003952569169960474
$ws.Range("O16").Value2 = 0.05533596837944664
$ws.Range("S16").Value2 = 0.08300395256916997
$ws.Range("F17").Value2 = 0.01463414634146342
$ws.Range("H17").Value2 = 0.1739837398373984
$ws.Range("I17").Value2 = 0.08130081300813008
$ws.Range("J17").Value2 = 0.4536585365853659
$ws.Range("K17").Value2 = 0.09268292682926829
$ws.Range("M17").Value2 = 0.01788617886178862
$ws.Range("O17").Value2 = 0.06178861788617886
$ws.Range("S17").Value2 = 0.1040650406504065
$ws.Range("F18").Value2 = 0.01798561151079137
$ws.Range("H18").Value2 = 0.197841726618705
$ws.Range("I18").Value2 = 0.08633093525179857
$ws.Range("J18").Value2 = 0.4892086330935252
$ws.Range("K18").Value2 = 0.07913669064748201
$ws.Range("M18").Value2 = 0.007194244604316547
$ws.Range("O18").Value2 = 0.05035971223021583
$ws.Range("S18").Value2 = 0.07194244604316546
$ws.Range("F19").Value2 = 0.01176470588235294
$ws.Range("H19").Value2 = 0.2007352941176471
$ws.Range("I19").Value2 = 0.06176470588235294
$ws.Range("J19").Value2 = 0.4036764705882353
$ws.Range("K19").Value2 = 0.09191176470588236
$ws.Range("M19").Value2 = 0.02573529411764706
$ws.Range("N19").Value2 = 0.001470588235294118
$ws.Range("O19").Value2 = 0.08602941176470588
$ws.Range("S19").Value2 = 0.1169117647058824
